$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.570003509521484
$ws.Range("B1").Value = 3.711336374282837
$ws.Range("C1").Value = 3.268029928207397
$ws.Range("D1").Value = 3.535899639129639
$ws.Range("E1").Value = 1.453612804412842
